$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C28").Value = 233
$ws.Range("D28").Value = 28
$ws.Range("E28").Value = 205
$ws.Range("F28").Value = 4.361370716510903
